$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1 (Suggestions/Issues row): turn the first trailing empty paragraph
# into a new paragraph of text.
# ---------------------------------------------------------------------------
$pD = $d.Paragraphs.Item(34)
$pD.Range.Text = "Tomcat configuration and learning curve. More research into this will be required."

# ---------------------------------------------------------------------------
# Step 2 (What's not going well and why row): replace the "n/a" (+ bookmark)
# paragraph with new wording, drop the old bookmark and remove one now
# redundant blank paragraph.
# ---------------------------------------------------------------------------
$pC = $d.Paragraphs.Item(28)
$pC.Range.InsertParagraphAfter()
$pCNew = $d.Paragraphs.Item(29)
$pCNew.Range.Text = "Migrating to Tomcat has proved to be harder than we thought."
# Delete the old "n/a" paragraph (this also removes the _GoBack bookmark that
# lived inside it).
$d.Paragraphs.Item(28).Range.Delete()
# Delete the now redundant blank paragraph that used to separate "n/a" from
# the bold empty paragraph.
$d.Paragraphs.Item(29).Range.Delete()

# ---------------------------------------------------------------------------
# Step 3 (Work to complete next reporting period row): drop one trailing
# blank paragraph.
# ---------------------------------------------------------------------------
$pB = $d.Paragraphs.Item(20)
$pB.Range.Delete()

# ---------------------------------------------------------------------------
# Step 4 (Work completed this reporting period row): drop one trailing blank
# paragraph and turn the last remaining one into the new home of the
# _GoBack bookmark.
# ---------------------------------------------------------------------------
$pA = $d.Paragraphs.Item(12)
$pA.Range.Delete()

$pABookmarkHost = $d.Paragraphs.Item(12)
# Bookmarks.Add mis-behaves on a truly collapsed (zero-length) range, so
# temporarily insert a placeholder character, bookmark that character, then
# remove the character again - the bookmark collapses down but stays put.
$pABookmarkHost.Range.InsertBefore("X")
$pABookmarkHost2 = $d.Paragraphs.Item(12)
$charRange = $d.Range($pABookmarkHost2.Range.Start, $pABookmarkHost2.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $charRange)
$d.Range($pABookmarkHost2.Range.Start, $pABookmarkHost2.Range.Start + 1).Text = ""
